$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the series-name cells under the benchmarking header row.
$ws.Range("C5").Value = "Sudoku"
$ws.Range("D5").Value = "Math Master"
$ws.Range("E5").Value = "Cálculo Mental"

# Move the active selection to E5 (matches the final cursor position in the
# authored workbook).
[void]$ws.Range("E5").Select()

# The category axis of the radar chart was flipped back to its natural
# (non-reversed) orientation.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$catAx = $chart.Axes(1)
$catAx.ReversePlotOrder = $false
